$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.7028063186220663
$ws.Cells.Item(2, 3).Value = 0.04515230317635144
$ws.Cells.Item(2, 4).Value = 0.1696043453175733
$ws.Cells.Item(2, 5).Value = 0.06506417358149008
$ws.Cells.Item(2, 6).Value = 2.882730028648723
$ws.Cells.Item(2, 9).Value = 2.379696603339184
$ws.Cells.Item(2, 11).Value = 0.5539231946166581
$ws.Cells.Item(2, 12).Value = 0.2350104351937787
$ws.Cells.Item(2, 13).Value = 0.2051997460670982
$ws.Cells.Item(3, 2).Value = 0.6830582863521499
$ws.Cells.Item(3, 3).Value = 0.03931581342145307
$ws.Cells.Item(3, 4).Value = 0.1687664912108247
$ws.Cells.Item(3, 5).Value = 0.06508989977272606
$ws.Cells.Item(3, 6).Value = 2.82979994787857
$ws.Cells.Item(3, 9).Value = 2.346872366106751
$ws.Cells.Item(3, 11).Value = 0.5277790340459489
$ws.Cells.Item(3, 12).Value = 0.2320207262118927
$ws.Cells.Item(3, 13).Value = 0.2009107545716837
$ws.Cells.Item(4, 2).Value = 0.6714924377058651
$ws.Cells.Item(4, 3).Value = 0.03572706087400945
$ws.Cells.Item(4, 4).Value = 0.1682345253070032
$ws.Cells.Item(4, 5).Value = 0.06511933386723268
$ws.Cells.Item(4, 6).Value = 2.798130741212773
$ws.Cells.Item(4, 9).Value = 2.327252651222693
$ws.Cells.Item(4, 11).Value = 0.5121768394279229
$ws.Cells.Item(4, 12).Value = 0.2302950901072478
$ws.Cells.Item(4, 13).Value = 0.1984151525753539
$ws.Cells.Item(5, 2).Value = 0.6669200062992218
$ws.Cells.Item(5, 3).Value = 0.03426322170970764
$ws.Cells.Item(5, 4).Value = 0.1680133123288385
$ws.Cells.Item(5, 5).Value = 0.06513476492832382
$ws.Cells.Item(5, 6).Value = 2.785433527912559
$ws.Cells.Item(5, 9).Value = 2.319391528250961
$ws.Cells.Item(5, 11).Value = 0.5059319503065467
$ws.Cells.Item(5, 12).Value = 0.2296195664454359
$ws.Cells.Item(5, 13).Value = 0.1974328564106287
$ws.Cells.Item(6, 2).Value = 0.6661692598584636
$ws.Cells.Item(6, 3).Value = 0.03402006512602895
$ws.Cells.Item(6, 4).Value = 0.1679763115718345
$ws.Cells.Item(6, 5).Value = 0.06513753500357122
$ws.Cells.Item(6, 6).Value = 2.783337728820001
$ws.Cells.Item(6, 9).Value = 2.318094282575231
$ws.Cells.Item(6, 11).Value = 0.5049018204343128
$ws.Cells.Item(6, 12).Value = 0.2295090688687083
$ws.Cells.Item(6, 13).Value = 0.1972718421135191
$ws.Cells.Item(7, 2).Value = 0.6714302023543723
$ws.Cells.Item(7, 3).Value = 0.03570732480650918
$ws.Cells.Item(7, 4).Value = 0.1682315599344086
$ws.Cells.Item(7, 5).Value = 0.06511952805242416
$ws.Cells.Item(7, 6).Value = 2.797958659629998
$ws.Cells.Item(7, 9).Value = 2.327146090978445
$ws.Cells.Item(7, 11).Value = 0.5120921607721129
$ws.Cells.Item(7, 12).Value = 0.2302858676374555
$ws.Cells.Item(7, 13).Value = 0.1984017645367651
$ws.Cells.Item(8, 2).Value = 0.6958810868926264
$ws.Cells.Item(8, 3).Value = 0.04314088021612861
$ws.Cells.Item(8, 4).Value = 0.1693190732269763
$ws.Cells.Item(8, 5).Value = 0.06507021750924213
$ws.Cells.Item(8, 6).Value = 2.864307142575498
$ws.Cells.Item(8, 9).Value = 2.368267657414165
$ws.Cells.Item(8, 11).Value = 0.544815111347873
$ws.Cells.Item(8, 12).Value = 0.2339567468040187
$ws.Cells.Item(8, 13).Value = 0.2036922974813429
$ws.Cells.Item(9, 2).Value = 0.7482707620662268
$ws.Cells.Item(9, 3).Value = 0.05768279880909688
$ws.Cells.Item(9, 4).Value = 0.1713139159674846
$ws.Cells.Item(9, 5).Value = 0.06508145503375573
$ws.Cells.Item(9, 6).Value = 3.001032914592429
$ws.Cells.Item(9, 9).Value = 2.45316824631027
$ws.Cells.Item(9, 11).Value = 0.6125691691670454
$ws.Cells.Item(9, 12).Value = 0.2420287271489627
$ws.Cells.Item(9, 13).Value = 0.2151610416290097
$ws.Cells.Item(10, 2).Value = 0.789478023103527
$ws.Cells.Item(10, 3).Value = 0.06835368576730616
$ws.Cells.Item(10, 4).Value = 0.1726973073643059
$ws.Cells.Item(10, 5).Value = 0.06515518555737287
$ws.Cells.Item(10, 6).Value = 3.105573187011828
$ws.Cells.Item(10, 9).Value = 2.51817819709602
$ws.Cells.Item(10, 11).Value = 0.6645544679181512
$ws.Cells.Item(10, 12).Value = 0.2484928594397644
$ws.Cells.Item(10, 13).Value = 0.2242557369543832
$ws.Cells.Item(11, 2).Value = 0.8088165435353005
$ws.Cells.Item(11, 3).Value = 0.07320719269779374
$ws.Cells.Item(11, 4).Value = 0.1733091725644584
$ws.Cells.Item(11, 5).Value = 0.06520287962823268
$ws.Cells.Item(11, 6).Value = 3.15403135467659
$ws.Cells.Item(11, 9).Value = 2.548332896438595
$ws.Cells.Item(11, 11).Value = 0.6886880951509227
$ws.Cells.Item(11, 12).Value = 0.2515497717200645
$ws.Cells.Item(11, 13).Value = 0.2285387979917175
$ws.Cells.Item(12, 2).Value = 0.816224913617873
$ws.Cells.Item(12, 3).Value = 0.07504511545374726
$ws.Cells.Item(12, 4).Value = 0.1735383896691332
$ws.Cells.Item(12, 5).Value = 0.06522296934332417
$ws.Cells.Item(12, 6).Value = 3.172511697295761
$ws.Cells.Item(12, 9).Value = 2.559835791264305
$ws.Cells.Item(12, 11).Value = 0.6978969355421896
$ws.Cells.Item(12, 12).Value = 0.2527240837257096
$ws.Cells.Item(12, 13).Value = 0.2301816658636326
$ws.Cells.Item(13, 2).Value = 0.8146255950982777
$ws.Cells.Item(13, 3).Value = 0.07464928391311787
$ws.Cells.Item(13, 4).Value = 0.1734891336856776
$ws.Cells.Item(13, 5).Value = 0.06521855252243292
$ws.Cells.Item(13, 6).Value = 3.168525821688974
$ws.Cells.Item(13, 9).Value = 2.557354696823211
$ws.Cells.Item(13, 11).Value = 0.6959105325818484
$ws.Cells.Item(13, 12).Value = 0.2524704308129344
$ws.Cells.Item(13, 13).Value = 0.2298269124990782
$ws.Cells.Item(14, 2).Value = 0.8094243251930493
$ws.Cells.Item(14, 3).Value = 0.07335839897328356
$ws.Cells.Item(14, 4).Value = 0.1733280800094192
$ws.Cells.Item(14, 5).Value = 0.06520449179268439
$ws.Cells.Item(14, 6).Value = 3.155549129347406
$ws.Cells.Item(14, 9).Value = 2.549277561631001
$ws.Cells.Item(14, 11).Value = 0.6894443088411037
$ws.Cells.Item(14, 12).Value = 0.2516460478925779
$ws.Cells.Item(14, 13).Value = 0.2286735375664435
$ws.Cells.Item(15, 2).Value = 0.8062495048719143
$ws.Cells.Item(15, 3).Value = 0.07256769861932355
$ws.Cells.Item(15, 4).Value = 0.1732291073886927
$ws.Cells.Item(15, 5).Value = 0.06519614324090028
$ws.Cells.Item(15, 6).Value = 3.147617513809962
$ws.Cells.Item(15, 9).Value = 2.544341028717128
$ws.Cells.Item(15, 11).Value = 0.6854926763597291
$ws.Cells.Item(15, 12).Value = 0.2511432675980245
$ws.Cells.Item(15, 13).Value = 0.2279697923104251
$ws.Cells.Item(16, 2).Value = 0.7882261405228235
$ws.Cells.Item(16, 3).Value = 0.06803649234598197
$ws.Cells.Item(16, 4).Value = 0.1726569716772772
$ws.Cells.Item(16, 5).Value = 0.06515235290805066
$ws.Cells.Item(16, 6).Value = 3.1024245339392
$ws.Cells.Item(16, 9).Value = 2.516219246147884
$ws.Cells.Item(16, 11).Value = 0.6629870584598336
$ws.Cells.Item(16, 12).Value = 0.2482954237890596
$ws.Cells.Item(16, 13).Value = 0.223978762960364
$ws.Cells.Item(17, 2).Value = 0.7773213036314814
$ws.Cells.Item(17, 3).Value = 0.06525663462836917
$ws.Cells.Item(17, 4).Value = 0.1723015349195975
$ws.Cells.Item(17, 5).Value = 0.06512910956836926
$ws.Cells.Item(17, 6).Value = 3.074931581649309
$ws.Cells.Item(17, 9).Value = 2.499116601288492
$ws.Cells.Item(17, 11).Value = 0.6493050024436968
$ws.Cells.Item(17, 12).Value = 0.2465781592442085
$ws.Cells.Item(17, 13).Value = 0.2215677439115709
$ws.Cells.Item(18, 2).Value = 0.7711049577074505
$ws.Cells.Item(18, 3).Value = 0.06365765834252102
$ws.Cells.Item(18, 4).Value = 0.1720954534197894
$ws.Cells.Item(18, 5).Value = 0.06511707319761229
$ws.Cells.Item(18, 6).Value = 3.059203256364697
$ws.Cells.Item(18, 9).Value = 2.489334309342482
$ws.Cells.Item(18, 11).Value = 0.6414811089141779
$ws.Cells.Item(18, 12).Value = 0.245601385379203
$ws.Cells.Item(18, 13).Value = 0.2201947167987583
$ws.Cells.Item(19, 2).Value = 0.769009798357672
$ws.Cells.Item(19, 3).Value = 0.0631162561025036
$ws.Cells.Item(19, 4).Value = 0.1720253946719197
$ws.Cells.Item(19, 5).Value = 0.06511322695397936
$ws.Cells.Item(19, 6).Value = 3.053892479840783
$ws.Cells.Item(19, 9).Value = 2.486031576457606
$ws.Cells.Item(19, 11).Value = 0.6388399103020674
$ws.Cells.Item(19, 12).Value = 0.245272547581834
$ws.Cells.Item(19, 13).Value = 0.2197321912809187
$ws.Cells.Item(20, 2).Value = 0.7784763641892027
$ws.Cells.Item(20, 3).Value = 0.06555256225841788
$ws.Cells.Item(20, 4).Value = 0.1723395416957914
$ws.Cells.Item(20, 5).Value = 0.06513144598982912
$ws.Cells.Item(20, 6).Value = 3.077849463305739
$ws.Cells.Item(20, 9).Value = 2.500931543460169
$ws.Cells.Item(20, 11).Value = 0.6507567534104339
$ws.Cells.Item(20, 12).Value = 0.2467598316074913
$ws.Cells.Item(20, 13).Value = 0.2218229803320355
$ws.Cells.Item(21, 2).Value = 0.8109497494041307
$ws.Cells.Item(21, 3).Value = 0.07373756215200444
$ws.Cells.Item(21, 4).Value = 0.1733754525602791
$ws.Cells.Item(21, 5).Value = 0.06520856675289899
$ws.Cells.Item(21, 6).Value = 3.159357158189493
$ws.Cells.Item(21, 9).Value = 2.551647730313789
$ws.Cells.Item(21, 11).Value = 0.6913416953439366
$ws.Cells.Item(21, 12).Value = 0.2518877352437556
$ws.Cells.Item(21, 13).Value = 0.2290117427810188
$ws.Cells.Item(22, 2).Value = 0.8326701364686642
$ws.Cells.Item(22, 3).Value = 0.07908701380317495
$ws.Cells.Item(22, 4).Value = 0.1740380193971731
$ws.Cells.Item(22, 5).Value = 0.06527079335998387
$ws.Cells.Item(22, 6).Value = 3.213386810096097
$ws.Cells.Item(22, 9).Value = 2.58528327228008
$ws.Cells.Item(22, 11).Value = 0.7182741672342843
$ws.Cells.Item(22, 12).Value = 0.2553366050600943
$ws.Cells.Item(22, 13).Value = 0.2338322165361717
$ws.Cells.Item(23, 2).Value = 0.8210320656962153
$ws.Cells.Item(23, 3).Value = 0.0762318666620132
$ws.Cells.Item(23, 4).Value = 0.1736857096924282
$ws.Cells.Item(23, 5).Value = 0.06523650203319242
$ws.Cells.Item(23, 6).Value = 3.184480479338276
$ws.Cells.Item(23, 9).Value = 2.567286431778186
$ws.Cells.Item(23, 11).Value = 0.7038624189371774
$ws.Cells.Item(23, 12).Value = 0.2534869595049827
$ws.Cells.Item(23, 13).Value = 0.2312482590411236
$ws.Cells.Item(24, 2).Value = 0.7779539962556044
$ws.Cells.Item(24, 3).Value = 0.06541877587562794
$ws.Cells.Item(24, 4).Value = 0.1723223642413174
$ws.Cells.Item(24, 5).Value = 0.06513038556140494
$ws.Cells.Item(24, 6).Value = 3.076530046893993
$ws.Cells.Item(24, 9).Value = 2.500110851613954
$ws.Cells.Item(24, 11).Value = 0.6501002856703906
$ws.Cells.Item(24, 12).Value = 0.246677664812708
$ws.Cells.Item(24, 13).Value = 0.2217075471421524
$ws.Cells.Item(25, 2).Value = 0.7336215740663476
$ws.Cells.Item(25, 3).Value = 0.05375173029669611
$ws.Cells.Item(25, 4).Value = 0.1707888236844504
$ws.Cells.Item(25, 5).Value = 0.06506689370159613
$ws.Cells.Item(25, 6).Value = 2.9633310307226
$ws.Cells.Item(25, 9).Value = 2.429740502188025
$ws.Cells.Item(25, 11).Value = 0.5938539922609891
$ws.Cells.Item(25, 12).Value = 0.2397514180240705
$ws.Cells.Item(25, 13).Value = 0.2119411528004917

Write-Output "Applied 216 cell updates (Case_5_29 380 kV pl_mw.xlsx)"
